$d = $word.ActiveDocument

# Locate the "Tema 2 Python" heading paragraph, then target the paragraph that
# immediately follows it (the pre-existing empty paragraph right before the
# section break) so the new content is inserted between the two.
$paras = $d.Paragraphs
$headingIndex = -1
for ($i = 1; $i -le $paras.Count; $i++) {
    if ($paras.Item($i).Range.Text -eq "Tema 2 Python`r") {
        $headingIndex = $i
    }
}
if ($headingIndex -eq -1) {
    throw "Could not find 'Tema 2 Python' paragraph"
}

$target = $paras.Item($headingIndex + 1)
$insertRange = $target.Range
$insertRange.Collapse(1)

$bodyFragment = '<w:p><w:pPr><w:rPr><w:bCs/></w:rPr></w:pPr><w:r><w:rPr><w:bCs/></w:rPr><w:br w:type="page"/></w:r></w:p><w:p><w:pPr><w:rPr><w:bCs/></w:rPr></w:pPr><w:r><w:rPr><w:bCs/></w:rPr><w:lastRenderedPageBreak/><w:t>19/10/2023</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:bCs/></w:rPr></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:bCs/></w:rPr><w:t>Computer</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:bCs/></w:rPr><w:t xml:space="preserve"> visión (visión por </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:bCs/></w:rPr><w:t>computardor</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:bCs/></w:rPr><w:t xml:space="preserve">) es la denominación técnica de utilizar IA (Deep </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:bCs/></w:rPr><w:t>learning</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:bCs/></w:rPr><w:t xml:space="preserve">) en </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:bCs/></w:rPr><w:t>imgenes</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:bCs/></w:rPr><w:t>.</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:bCs/></w:rPr></w:pPr><w:r><w:rPr><w:bCs/></w:rPr><w:t xml:space="preserve">Diferentes tipos de problemas que </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:bCs/></w:rPr><w:t>see</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:bCs/></w:rPr><w:t xml:space="preserve"> pueden encontrar</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="3"/></w:numPr><w:rPr><w:bCs/></w:rPr></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:bCs/></w:rPr><w:t>Classificación</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:bCs/></w:rPr><w:t xml:space="preserve"> / Regresión</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:bCs/></w:rPr></w:pPr></w:p><w:p><w:pPr><w:rPr><w:bCs/></w:rPr></w:pPr><w:r><w:rPr><w:bCs/></w:rPr><w:t>Introducción IA en imagen médica</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:bCs/></w:rPr></w:pPr><w:r><w:rPr><w:bCs/></w:rPr><w:t xml:space="preserve">Detección vs segmentación: señalar </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:bCs/></w:rPr><w:t>uan</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:bCs/></w:rPr><w:t xml:space="preserve"> región vs identificar cada uno de los píxeles</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:bCs/></w:rPr></w:pPr><w:r><w:rPr><w:bCs/></w:rPr><w:t>La segmentación es de gran importancia en imagen médica</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:bCs/></w:rPr></w:pPr></w:p><w:p><w:pPr><w:rPr><w:bCs/></w:rPr></w:pPr><w:r><w:rPr><w:bCs/></w:rPr><w:t>Entrenamiento:</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:bCs/></w:rPr></w:pPr><w:r><w:rPr><w:bCs/></w:rPr><w:t>Paso 1: Preparación de la base de datos</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="3"/></w:numPr><w:rPr><w:bCs/></w:rPr></w:pPr><w:r><w:rPr><w:bCs/></w:rPr><w:t>Entrenamiento</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="3"/></w:numPr><w:rPr><w:bCs/></w:rPr></w:pPr><w:r><w:rPr><w:bCs/></w:rPr><w:t xml:space="preserve">Validación: ajustar </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:bCs/></w:rPr><w:t>hyperparametros</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:bCs/></w:rPr><w:t>. Ejemplo: numero de capas y neuronas</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="3"/></w:numPr><w:rPr><w:bCs/></w:rPr></w:pPr><w:r><w:rPr><w:bCs/></w:rPr><w:t>Test</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:bCs/></w:rPr></w:pPr><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:bCs/></w:rPr><w:t>Concepto covariables</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:bCs/></w:rPr><w:t>:</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:bCs/></w:rPr></w:pPr><w:r><w:rPr><w:bCs/></w:rPr><w:t>Estratificación: como se hacen las particiones. Garantizar que el numero de pruebas en cada partición es el mismo. Ejemplo: el mismo % de hombre y mujeres en entrenamiento y test</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:bCs/></w:rPr></w:pPr></w:p><w:p><w:pPr><w:rPr><w:bCs/></w:rPr></w:pPr><w:r><w:rPr><w:bCs/></w:rPr><w:t xml:space="preserve">Paso 2: ajuste de </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:bCs/></w:rPr><w:t>hiperparametros</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p><w:p><w:pPr><w:rPr><w:bCs/></w:rPr></w:pPr><w:r><w:rPr><w:bCs/></w:rPr><w:t>Parámetros que se han de seleccionar de manera manual antes de empezar el aprendizaje</w:t></w:r></w:p>'

$packageXml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>' + $bodyFragment + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$insertRange.InsertXML($packageXml)

Write-Output "Inserted new content after 'Tema 2 Python' heading."
